$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170-245 down to 171-246.
$ws.Rows("170").Insert()

# Populate the newly inserted row 170 with the new weekly record.
$ws.Range("A170").Value = 8
$ws.Range("B170").Value = "Terminal La Palmera de La Serena"
$ws.Range("C170").Value = "Coquimbo"
$ws.Range("D170").Value = 44839
$ws.Range("D170").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E170").Value = 4
$ws.Range("F170").Value = 100112037
$ws.Range("G170").Value = "Cebollín"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 1400
$ws.Range("K170").Value = 1400
$ws.Range("L170").Value = 1600
$ws.Range("M170").Value = 1500
$ws.Range("N170").Value = "$/paquete 6 unidades"
$ws.Range("O170").Value = "Provincia del Elquí"
$ws.Range("P170").Value = 250
$ws.Range("Q170").Value = 6
$ws.Range("R170").Value = "Hortaliza"
